# "resolve by type with depends"
#
# 1. Iterations sheet (sheet1): add a new "Register with simple parameters"
#    line under Iteration 2 (row 7), highlight the Iteration 2 header with
#    the same green fill used by the Iteration 1 header, extend the green
#    B-column highlight strip down through the new row, underline the
#    "LifeStyle Transient" line, insert a new "LifeStyle Singleton" line
#    under the LifeStyle group (pushing Iteration 5 down a row), and move
#    the selection.
# 2. Add a new "Limitations" sheet after "Iterations" documenting a caveat.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Iterations")

# --- Iteration 2 header: pick up the same bold+green style as the
#     Iteration 1 header (A1) so the two section headers match. ---
$ws.Range("A1").Copy()
$ws.Range("A5").PasteSpecial(-4122)

# --- extend the green B-column accent strip down to the new row ---
$ws.Range("B2").Copy()
$ws.Range("B5:B7").PasteSpecial(-4122)

# --- new bullet under "Iteration 2" ---
$ws.Range("A7").Value = "Register with simple parameters"

# --- underline the "LifeStyle Transient" bullet ---
$ws.Range("A13").Font.Underline = $true

# --- make room for a new "LifeStyle Singleton" bullet: push the
#     "Iteration 5" header (and everything below it) down one row ---
$ws.Rows("16").Insert()
$ws.Range("A15").Value = "LifeStyle Singleton"

$ws.Range("B10").Select()

# --- new "Limitations" sheet, placed right after "Iterations" ---
$wsLimitations = $wb.Worksheets.Add($null, $ws)
$wsLimitations.Name = "Limitations"
$wsLimitations.Range("A1").Value = "Services with only one constructor supported"
$wsLimitations.Columns("A").ColumnWidth = 54.109375
$wsLimitations.Range("A3").Select()

$ws.Select()
